# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.883.06"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.602.14"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("D6").Value = "154.95"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "'6.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "3.059.09"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "60.909.92"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "21.68"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "2.607.97"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "355.78"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").Value = "10.59"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "61.03"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "0.426"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "2.716.54"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "0.989"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +9.77%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "149.21"
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "0.917"
$ws.Range("E37").Value = "  +7.49%  "
$ws.Range("D38").Value = "'0.900"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "3.81"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "36.48"
$ws.Range("D42").Value = "291.84"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D44").Value = "0.623"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "19.64"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").Value = "4.92"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "0.0238"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").Value = "10.34"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").Value = "19.21"
$ws.Range("E51").Value = "  +8.35%  "
